$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.002.31"
$ws.Range("E2").Value = "  +0.87%  "

$ws.Range("D3").Value = "2.393.91"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "507.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.66%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("E8").Value = "  -0.27%  "

$ws.Range("D9").Value = "2.398.95"
$ws.Range("E9").Value = "  -0.31%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0987"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.41%  "

$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("E12").Value = "  +4.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.62%  "

$ws.Range("D14").Value = "2.817.79"
$ws.Range("E14").Value = "  +0.26%  "

$ws.Range("D15").Value = "56.953.39"
$ws.Range("E15").Value = "  +1.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.80%  "

$ws.Range("E17").Value = "  +1.21%  "

$ws.Range("D18").Value = "2.406.29"
$ws.Range("E18").Value = "  +0.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "310.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.19%  "

$ws.Range("E22").Value = "  -1.58%  "

$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.373"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.151"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.31%  "

$ws.Range("D30").Value = "0.0₃0729"
$ws.Range("E30").Value = "  +0.97%  "

$ws.Range("E31").Value = "  -0.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.42%  "

$ws.Range("E37").Value = "  -1.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.824"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.66%  "

$ws.Range("E41").Value = "  -0.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "131.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.70%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.47%  "

$ws.Range("E45").Value = "  +0.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0912"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "249.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0486"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.78%  "

$ws.Range("E49").Value = "  +0.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.25%  "

$ws.Range("E51").Value = "  +2.03%  "
